# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.096.57"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.569.66"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'314.89"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("D10").Value = "'35.47"
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'7.45"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("D13").Value = "2.964.08"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").Value = "2.595.93"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("D16").Value = "'15.06"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "43.132.03"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  +2.71%  "
$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = "  -3.44%  "
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").Value = "'69.31"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").Value = "'253.40"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'2.97"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "'26.90"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "'40.20"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").Value = "'10.27"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  -4.06%  "
$ws.Range("D32").Value = "'154.83"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'18.98"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  +6.94%  "
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("E41").Value = "  -5.56%  "
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("D43").Value = "'0.0305"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'3.24"
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("D46").Value = "1.999.20"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").Value = "2.815.01"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").Value = "'82.84"
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("D50").Value = "'74.73"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").Value = "'0.194"
$ws.Range("E51").Value = "  +2.38%  "
